# Rename the worksheet tab (drops the internal plate identifier from the name)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Example plate"

# Update the active selection / view to B32 (also clears the scrolled
# "topLeftCell" view state that pointed at A28/C29)
$ws.Range("B32").Select()
